$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '308.96'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = '-5.71%'
$ws.Cells.Item(2,5).Style = "Normal"

# Row 3
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '40.35'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = '-9.03%'
$ws.Cells.Item(3,5).Style = "Normal"

# Row 4
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '5.052'
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = '-4.54%'
$ws.Cells.Item(4,5).Style = "Normal"

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '0.07786'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = '-7.00%'
$ws.Cells.Item(5,5).Style = "Normal"

# Row 6
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = '-1.94%'
$ws.Cells.Item(6,5).Style = "Normal"

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '1.663'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = '-14.06%'
$ws.Cells.Item(7,5).Style = "Normal"

# Row 8
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.9095'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = '-6.29%'
$ws.Cells.Item(8,5).Style = "Normal"

# Row 9
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.1039'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = '-8.49%'
$ws.Cells.Item(9,5).Style = "Normal"

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.1755'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = '-7.67%'
$ws.Cells.Item(10,5).Style = "Normal"

# Row 11
$ws.Cells.Item(11,2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11,3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '0.09003'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = '-6.79%'
$ws.Cells.Item(11,5).Style = "Normal"

# Row 12
$ws.Cells.Item(12,2).Value = 'BitrueCoin'
$ws.Cells.Item(12,3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = '0.04429'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = '-3.71%'
$ws.Cells.Item(12,5).Style = "Normal"

# Row 13
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '7.132'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = '-16.18%'
$ws.Cells.Item(13,5).Style = "Normal"

# Row 14
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '0.1057'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = '-0.42%'
$ws.Cells.Item(14,5).Style = "Normal"

# Row 15
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.001250'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = '-3.67%'
$ws.Cells.Item(15,5).Style = "Normal"

# Row 16
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '0.005840'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = '0.83%'
$ws.Cells.Item(16,5).Style = "Normal"

# Row 17
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '3.362'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = '-0.76%'
$ws.Cells.Item(17,5).Style = "Normal"

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '0.3367'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = '0.29%'
$ws.Cells.Item(19,5).Style = "Normal"

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '0.1390'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = '0.15%'
$ws.Cells.Item(20,5).Style = "Normal"

# Row 21
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = '0.2853'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = '10.74%'
$ws.Cells.Item(21,5).Style = "Normal"

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '0.04172'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = '0.29%'
$ws.Cells.Item(22,5).Style = "Normal"

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '0.001217'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value = '-1.29%'
$ws.Cells.Item(23,5).Style = "Normal"

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '0.004103'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value = '-6.95%'
$ws.Cells.Item(24,5).Style = "Normal"

# Row 25
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = '0.0001226'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value = '-5.71%'
$ws.Cells.Item(25,5).Style = "Normal"

# Row 26
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '0.0002994'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value = '0.52%'
$ws.Cells.Item(26,5).Style = "Normal"

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '0.02410'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).NumberFormat = "@"
$ws.Cells.Item(38,5).Value = '-11.23%'
$ws.Cells.Item(38,5).Style = "Normal"

# Row 39
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '0.05214'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value = '-7.11%'
$ws.Cells.Item(39,5).Style = "Normal"

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '0.007971'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value = '1.93%'
$ws.Cells.Item(40,5).Style = "Normal"

# Row 41
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.1330'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value = '-5.96%'
$ws.Cells.Item(41,5).Style = "Normal"

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '0.007577'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value = '3.79%'
$ws.Cells.Item(42,5).Style = "Normal"

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.001990'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value = '-2.92%'
$ws.Cells.Item(43,5).Style = "Normal"

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '0.008049'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value = '-7.17%'
$ws.Cells.Item(44,5).Style = "Normal"

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '0.3361'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value = '-4.37%'
$ws.Cells.Item(45,5).Style = "Normal"

# Row 46
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '0.00006740'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value = '-2.47%'
$ws.Cells.Item(46,5).Style = "Normal"

# Row 47
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '0.00000000754'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value = '0.50%'
$ws.Cells.Item(47,5).Style = "Normal"

# Row 48
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.003322'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).NumberFormat = "@"
$ws.Cells.Item(48,5).Value = '-4.89%'
$ws.Cells.Item(48,5).Style = "Normal"

# Row 49
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '0.004122'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).NumberFormat = "@"
$ws.Cells.Item(49,5).Value = '16.75%'
$ws.Cells.Item(49,5).Style = "Normal"

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '0.00002111'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).NumberFormat = "@"
$ws.Cells.Item(50,5).Value = '0.50%'
$ws.Cells.Item(50,5).Style = "Normal"

# Row 51
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '0.0002010'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).NumberFormat = "@"
$ws.Cells.Item(51,5).Value = '0.50%'
$ws.Cells.Item(51,5).Style = "Normal"
